# Re-upload / re-save of data/type5.xlsx:
#  - A2 and A3 change from 8 to 5
#  - the active selection moves from A4:C5 to the single cell A7
#  - the window size recorded in the file (bookViews) also changed, but
#    that is a host-application chrome setting (not part of the workbook
#    content/object model) so there is nothing to set through COM for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two data cells: A2 and A3 from 8 to 5
$ws.Range("A2").Value = 5
$ws.Range("A3").Value = 5

# Move/collapse the selection to A7 (activeCell="A7", sqref="A7")
$ws.Range("A7").Select()
